$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.6606524410359556
$ws.Range("C2").Value = 10.34677158129881
$ws.Range("D2").Value = 401567.231247708
$ws.Range("E2").Value = 91228006295.30009
$ws.Range("G2").Value = 91228407873.53877
